# The data rows (2-23) of the sheet were re-sorted/shuffled upstream; apply
# that shuffle here as a permutation of whole rows (same columns, just rows
# swapped around). Rows 12, 15 and 18 stay exactly where they are.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# new row -> source row (row numbers as they exist in the "before" sheet)
$rowMap = @{
    2  = 10
    3  = 5
    4  = 21
    5  = 13
    6  = 2
    7  = 9
    8  = 16
    9  = 17
    10 = 6
    11 = 8
    12 = 12
    13 = 7
    14 = 23
    15 = 15
    16 = 14
    17 = 22
    18 = 18
    19 = 3
    20 = 4
    21 = 11
    22 = 19
    23 = 20
}

# Snapshot every source row (A:R) before any writes happen, since several
# rows both give and receive data (the permutation has cycles).
$firstRow = 2
$lastRow = 23
$snapshot = $ws.Range("A$firstRow`:R$lastRow").Value2

$numCols = 18

for ($newRow = $firstRow; $newRow -le $lastRow; $newRow++) {
    $oldRow = $rowMap[$newRow]
    if ($oldRow -eq $newRow) {
        continue
    }
    $srcOffset = $oldRow - $firstRow + 1
    $rowValues = New-Object 'object[,]' 1,$numCols
    for ($col = 1; $col -le $numCols; $col++) {
        $rowValues[0,$col-1] = $snapshot[$srcOffset,$col]
    }
    $ws.Range("A$newRow`:R$newRow").Value = $rowValues
}
